# Update test data for 5,24,40V,BatteryStandby and AC Calculations test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# C8: CPU 800 -> CPU 801
$ws.Range("C8").Value = "CPU 801"

# C9: CPU 800 -> "" (empty text, matches style/quotePrefix already used by C10/C11)
$ws.Range("C10").Copy()
$ws.Range("C9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C9").Value = "'"

# A11: Pro32xD -> Pro885D
$ws.Range("A11").Value = "Pro885D"

# Remove rows 12 and 13 (MX2-100 / P885D entries no longer needed)
$ws.Rows("12:13").Delete()

# Update the active selection to match the saved view state
$ws.Range("C9").Select()
